$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "strain" values (column F) for rows 2-22.
# New shared strings introduced, in first-use order:
#   KN99alpha, TDY1379, TDY1375, TDY1366, TDY1378, TDY1374
$ws.Range("F2").Value = "KN99alpha"
$ws.Range("F3").Value = "KN99alpha"
$ws.Range("F4").Value = "KN99alpha"

$ws.Range("F5").Value = "TDY1379"
$ws.Range("F6").Value = "TDY1379"
$ws.Range("F7").Value = "TDY1379"

$ws.Range("F8").Value = "TDY1375"
$ws.Range("F9").Value = "TDY1375"
$ws.Range("F10").Value = "TDY1375"

$ws.Range("F11").Value = "TDY1366"
$ws.Range("F12").Value = "TDY1366"
$ws.Range("F13").Value = "TDY1366"

$ws.Range("F14").Value = "KN99alpha"
$ws.Range("F15").Value = "KN99alpha"
$ws.Range("F16").Value = "KN99alpha"

$ws.Range("F17").Value = "TDY1378"
$ws.Range("F18").Value = "TDY1378"
$ws.Range("F19").Value = "TDY1378"

$ws.Range("F20").Value = "TDY1374"
$ws.Range("F21").Value = "TDY1374"
$ws.Range("F22").Value = "TDY1374"

# Update the sheet's selection / active cell to match the saved view state.
$ws.Range("F21:F22").Select() | Out-Null
